$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTP2021")

# Update row 2: directory changes from 2005_TM152_IPA_02 to 2005_TM152_IPA_03
$ws.Range("B2").Value = "2005_TM152_IPA_03"

# Update row 3: year 2035 -> 2015, directory -> 2015_TM152_IPA_16, category aoc1421 -> IP
$ws.Range("A3").Value = 2015
$ws.Range("B3").Value = "2015_TM152_IPA_16"
$ws.Range("C3").Value = "IP"

# Update row 4: year stays 2035, directory -> 2035_TM152_IPA_00, category aoc1562 -> IP
$ws.Range("A4").Value = 2035
$ws.Range("B4").Value = "2035_TM152_IPA_00"
$ws.Range("C4").Value = "IP"

# Remove row 5 entirely (previously 2035 / 2035_TM152_IPA_aoc1795_00 / IP_aoc1795)
$ws.Range("A5:C5").EntireRow.Delete()

# Update selection to match the new last-used cell
$ws.Range("A4").Select()
